# Sync attendance_reports: swap order of "Recorded By" entries in column G
# so that entries formatted as "dnasr281@gmail.com, X" become "X, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -ne $null -and $val -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, dnasr281@gmail.com"
    }
}
